# Risk assessment workbook update
# - Renames Sheet1..Sheet4 to their descriptive names
# - Adds two new sheets ("Track Contoller", "Track Model") for the remaining
#   subsystems, with "Track Contoller" carrying a "T" placeholder in A1
# - Restores per-sheet row heights that were nudged by Excel's autofit
# - Re-applies the saved selection/scroll state for each sheet, leaving
#   "CTC Office" (sheet 2) as the active tab, matching the saved file

$wb = $excel.ActiveWorkbook

# --- Rename existing sheets -------------------------------------------------
$wb.Worksheets.Item(1).Name = "Project RA"
$wb.Worksheets.Item(2).Name = "CTC Office"
$wb.Worksheets.Item(3).Name = "Train Controller"
$wb.Worksheets.Item(4).Name = "Train Model"

$wsProjectRA      = $wb.Worksheets.Item(1)
$wsCtcOffice      = $wb.Worksheets.Item(2)
$wsTrainCtrl      = $wb.Worksheets.Item(3)
$wsTrainModel     = $wb.Worksheets.Item(4)

# --- Add the two new subsystem sheets ---------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTrackCtrl = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet)
$wsTrackCtrl.Name = "Track Contoller"
$wsTrackCtrl.Range("A1").Value = "T"

$lastSheet2 = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsTrackModel = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $lastSheet2)
$wsTrackModel.Name = "Track Model"

# --- Row height touch-ups ----------------------------------------------------
$wsProjectRA.Rows.Item(4).RowHeight = 30.6
$wsProjectRA.Rows.Item(6).RowHeight = 105
$wsProjectRA.Rows.Item(7).RowHeight = 60
$wsProjectRA.Rows.Item(8).RowHeight = 45.95

$wsCtcOffice.Rows.Item(3).RowHeight = 62.1
$wsCtcOffice.Rows.Item(4).RowHeight = 87.6
$wsCtcOffice.Rows.Item(5).RowHeight = 134.45
$wsCtcOffice.Rows.Item(7).RowHeight = 90
$wsCtcOffice.Rows.Item(8).RowHeight = 90
$wsCtcOffice.Rows.Item(9).RowHeight = 150
$wsCtcOffice.Rows.Item(11).RowHeight = 165

# --- Selections / scroll position for the sheets whose view state changed ---
# (Train Controller / Train Model stay untouched - no selection recorded there)
$wsProjectRA.Activate()
$null = $wsProjectRA.Range("D13").Select()

$wsTrackCtrl.Activate()
$null = $wsTrackCtrl.Range("A1").Select()

$wsTrackModel.Activate()
$null = $wsTrackModel.Range("D24").Select()

# CTC Office is the sheet that was left active/selected when the workbook was
# last saved, so activate it last.
$wsCtcOffice.Activate()
$null = $wsCtcOffice.Range("H8").Select()
